$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Resize the existing table ("표1") to its new extent first. Doing this
#    BEFORE the header text is (re)written lets the engine resync each
#    ListColumn's name against its header cell text when the workbook is
#    saved.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$null = $lo.Resize($ws.Range("A1:F15"))

# ---------------------------------------------------------------------------
# 2. Header row - insert "Name" after "Id", and "ItemType"/"EngineId" before
#    "ItemIcon".
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "Id"
$ws.Cells.Item(1,2).Value = "Name"
$ws.Cells.Item(1,3).Value = "DisplayName"
$ws.Cells.Item(1,4).Value = "ItemType"
$ws.Cells.Item(1,5).Value = "EngineId"
$ws.Cells.Item(1,6).Value = "ItemIcon"

# ---------------------------------------------------------------------------
# 3. Existing resource rows (2-9): shift the old DisplayName out of column B
#    into column C, add the internal "Name" key in column B and a new
#    "ItemType" value ("Resource") in column D. Columns E/F stay untouched.
# ---------------------------------------------------------------------------
$resourceRows = @(
  @(1001, "Stone_Material",  "돌 광물"),
  @(1002, "Stone_Resource",  "돌 자재"),
  @(1003, "Copper_Material", "구리 광물"),
  @(1004, "Copper_Resource", "구리 자재"),
  @(1005, "Iron_Material",   "철 광물"),
  @(1006, "Iron_Resource",   "철 자재"),
  @(1007, "Gold_Material",   "금 광물"),
  @(1008, "Gold_Resource",   "금 자재")
)

$r = 2
foreach ($row in $resourceRows) {
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = "Resource"
  $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4. Brand new "Engine" rows (10-15): Id, Name, DisplayName, ItemType and the
#    numeric EngineId. Column F (ItemIcon) is left empty for all of them.
# ---------------------------------------------------------------------------
$engineRows = @(
  @(1101, "GreenEngine_1", "초록 엔진 Lv.1", 201001),
  @(1102, "GreenEngine_2", "초록 엔진 Lv.2", 201002),
  @(1103, "GreenEngine_3", "초록 엔진 Lv.3", 201003),
  @(1104, "RedEngine_1",   "빨강 엔진 Lv.1", 202001),
  @(1105, "RedEngine_2",   "빨강 엔진 Lv.2", 202002),
  @(1106, "RedEngine_3",   "빨강 엔진 Lv.3", 202003)
)

$r = 10
foreach ($row in $engineRows) {
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = "Engine"
  $ws.Cells.Item($r,5).Value = $row[3]
  $r = $r + 1
}

# ---------------------------------------------------------------------------
# 5. Column widths for the two freshly-introduced columns (D, E) match the
#    width already used by column B.
# ---------------------------------------------------------------------------
$existingWidth = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = $existingWidth
$ws.Columns.Item(5).ColumnWidth = $existingWidth

# ---------------------------------------------------------------------------
# 6. Move the active selection to match where the author ended up.
# ---------------------------------------------------------------------------
$null = $ws.Range("G8:G9").Select()
